$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, shifting all existing data rows (17-115) down to (18-116).
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44547
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112045
$ws.Cells.Item(17, 7).Value = "Zapallo"
$ws.Cells.Item(17, 8).Value = "Camote"
$ws.Cells.Item(17, 9).Value = "1a nueva(o)"
$ws.Cells.Item(17, 10).Value = 600
$ws.Cells.Item(17, 11).Value = 600
$ws.Cells.Item(17, 12).Value = 650
$ws.Cells.Item(17, 13).Value = 625
$ws.Cells.Item(17, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(17, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 16).Value = 625
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
